# resetPassword test and some adjustments
$wb = $excel.ActiveWorkbook

# --- Parameters sheet (new appium/device config rows) ---
$params = $wb.Worksheets.Item("Parameters")

$params.Cells.Item(14, 1).Value = "appiumVersion"
$params.Cells.Item(14, 2).Value = "1.6.4"

$params.Cells.Item(15, 1).Value = "deviceName"
$params.Cells.Item(15, 2).Value = "iPhone 7 Simulator"

$params.Cells.Item(16, 1).Value = "deviceOrientation"
$params.Cells.Item(16, 2).Value = "portrait"

$params.Cells.Item(17, 1).Value = "platformVersion"
$params.Cells.Item(17, 2).Value = "10.3"

$params.Cells.Item(18, 1).Value = "platformName"
$params.Cells.Item(18, 2).Value = "iOS"

$params.Range("A14:B18").Select()

# --- Scenarios sheet ---
$scenarios = $wb.Worksheets.Item("Scenarios")

# Flip Execute column to "Y" for the existing rows 2-12 (tests that were N)
for ($r = 2; $r -le 12; $r++) {
    $scenarios.Cells.Item($r, 1).Value = "Y"
}

# Insert a new row 13 for the PasswordTest/resetPassword case, pushing
# the ops.* rows down to 14 and 15.
$scenarios.Rows.Item(13).Insert()

$scenarios.Cells.Item(13, 1).Value = "Y"
$scenarios.Cells.Item(13, 2).Value = "patient.tests.PasswordTest"
$scenarios.Cells.Item(13, 4).Value = "resetPassword"

$scenarios.Range("A15").Select()

$wb.Save()
